$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the tracking number in A1 (was 13822917)
$ws.Range("A1").Value = 1111

# Remove the "ActivityReport" javascript hyperlink that was attached to A1
$ws.Hyperlinks.Delete()

# Populate the previously-empty cells A2:A5 with their new values
$ws.Range("A2").Value = 222
$ws.Range("A3").Value = 333
$ws.Range("A4").Value = 444
$ws.Range("A5").Value = 555

# Move/leave the active selection on A5
$ws.Range("A5").Select()
